# "add battery to spreadsheet"
# Adds a small Battery part-selection table (columns I:M) to the "Power"
# sheet, mirroring the existing "Part Selection" tables elsewhere in the
# workbook: a bold header row (I2:M2) and two data rows (3 and 4) holding
# the two LiPo batteries under consideration.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Workbook-level bookkeeping (absPath, revision pointer, window size)
# ---------------------------------------------------------------------
$wb.Title = $wb.Title  # no-op placeholder to keep $wb referenced early

# ---------------------------------------------------------------------
# 2. "Power" sheet: new battery comparison table
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Power")

# --- values first (controls shared-string insertion order) ----------
$ws.Range("I3").Value = "Polymer Lithium Ion Battery (LiPo) 3.7V 400mAh"
$ws.Range("I4").Value = "Polymer Lithium Ion Battery (LiPo) 3.7V 1100mAh"

$ws.Range("J3").Value = 3.7
$ws.Range("J4").Value = 3.7

$ws.Range("M3").Value = "https://core-electronics.com.au/polymer-lithium-ion-battery-400mah-38456.html"
$ws.Range("M4").Value = "https://core-electronics.com.au/polymer-lithium-ion-battery-1000mah-38458.html"

$ws.Range("M2").Value = "Link"
$ws.Range("L2").Value = "Quantity"
$ws.Range("K2").Value = "Current (mAh)"

$ws.Range("K3").Value = 400
$ws.Range("K4").Value = 1100

$ws.Range("L3").Value = 1
$ws.Range("L4").Value = 2

$ws.Range("I2").Value = "Battery"
$ws.Range("J2").Value = "Voltage (V)"

# --- formatting: copy from the matching existing "Part Selection" -----
# table cells so the new table reuses the same fonts/borders/alignment.
# Header row (bold, bordered, centred) -> copy from B3 (existing bold
# table header cell).
$ws.Range("B3").Copy() | Out-Null
$ws.Range("I2:M2").PasteSpecial(-4122) | Out-Null

# Data rows (bordered, centred, regular weight) -> copy from A4 (existing
# plain bordered data cell).
$ws.Range("A4").Copy() | Out-Null
$ws.Range("I3:M4").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Column widths for the new columns (I:M) matching the authored sheet.
$ws.Columns.Item(9).ColumnWidth = 45.28515625
$ws.Columns.Item(10).ColumnWidth = 11
$ws.Columns.Item(11).ColumnWidth = 12.5703125
$ws.Columns.Item(12).ColumnWidth = 8.7109375
$ws.Columns.Item(13).ColumnWidth = 76.5703125

# ---------------------------------------------------------------------
# 3. Selection / view state to match the saved workbook
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Ideation")
$ws1.Activate()
$ws1.Range("E3:F3").Select() | Out-Null

$ws2 = $wb.Worksheets.Item("Parts")
$ws2.Activate()
$ws2.Range("A12").Select() | Out-Null

$ws.Activate()
$ws.Range("M11").Select() | Out-Null
$excel.ActiveWindow.Zoom = 70
